$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2..22, replacing the old Strike# based values.
$newValues = @{
    2  = 2
    3  = 6
    4  = 2
    5  = 3
    6  = 6
    7  = 9
    8  = 1
    9  = 6
    10 = 2
    11 = 2
    12 = 3
    13 = 4
    14 = 3
    15 = 5
    16 = 1
    17 = 4
    18 = 2
    19 = 5
    20 = 1
    21 = 2
    22 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
